# Update Handback report timestamps (commit: "Generate Report for Handback")
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first file row
$wsOverview.Range("G2").Value = "2016-08-23 01:01:24"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn.Range("H2").Value = "2016-08-23 01:01:19"
$wsZhCn.Range("K2").Value = "2016-08-23 01:01:44"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe.Range("H2").Value = "2016-08-23 01:01:24"
$wsDeDe.Range("K2").Value = "2016-08-23 01:01:51"
